# reversal report, tabbed, table - norms
# Rebuild the raw->ss lookup tables on sheets "5.0","5.6","6.0" (continuous raw
# scores 1-30, replacing the old banded/inline-string rows) and patch the
# handful of standard-score cells that shifted on sheets "6.6","7.0","7.6","8.0".

$wb = $excel.ActiveWorkbook

function Set-LookupColumn($ws, $Values) {
    # $Values holds B2..B31 in order (paired with raw scores 1..30 in col A)
    $n = $Values.Length
    $arr = New-Object 'object[,]' $n,2
    for ($i = 0; $i -lt $n; $i++) {
        $arr[$i,0] = $i + 1
        $arr[$i,1] = $Values[$i]
    }
    $rng = $ws.Range("A2:B$($n + 1)")
    $rng.Value = $arr
}

# ---- Sheet "5.0" (sheet1) : A1:B23 -> A1:B31, raw 1..30 ----
$ws1 = $wb.Worksheets.Item(1)
Set-LookupColumn $ws1 @(
    82,87,87,94,97,99,101,103,105,107,
    109,110,112,114,115,117,119,121,123,124,
    126,128,130,130,130,130,130,130,130,130
)

# ---- Sheet "5.6" (sheet2) : A1:B29 -> A1:B31, raw 1..30 ----
$ws2 = $wb.Worksheets.Item(2)
Set-LookupColumn $ws2 @(
    75,82,82,82,89,91,93,95,97,98,
    100,102,103,105,106,108,110,111,113,115,
    116,118,120,122,124,126,129,130,130,130
)

# ---- Sheet "6.0" (sheet3) : A1:B25 -> A1:B31, raw 1..30 ----
$ws3 = $wb.Worksheets.Item(3)
Set-LookupColumn $ws3 @(
    73,77,77,77,77,77,88,90,91,93,
    94,96,97,99,100,102,103,105,107,108,
    110,112,114,116,118,120,123,126,130,130
)

# ---- Sheet "6.6" (sheet4): targeted SS corrections, raw row -> new value ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value  = 59
$ws4.Range("B6").Value  = 78
$ws4.Range("B8").Value  = 81
$ws4.Range("B30").Value = 120

# ---- Sheet "7.0" (sheet5): targeted SS corrections ----
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value  = 67
$ws5.Range("B3").Value  = 69
$ws5.Range("B4").Value  = 71
$ws5.Range("B6").Value  = 74
$ws5.Range("B8").Value  = 77
$ws5.Range("B9").Value  = 78
$ws5.Range("B12").Value = 82
$ws5.Range("B15").Value = 86
$ws5.Range("B18").Value = 90
$ws5.Range("B20").Value = 93
$ws5.Range("B30").Value = 118

# ---- Sheet "7.6" (sheet6): targeted SS corrections ----
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B2").Value  = 65
$ws6.Range("B3").Value  = 66
$ws6.Range("B4").Value  = 68
$ws6.Range("B5").Value  = 69
$ws6.Range("B6").Value  = 71
$ws6.Range("B7").Value  = 72
$ws6.Range("B8").Value  = 73
$ws6.Range("B9").Value  = 75
$ws6.Range("B10").Value = 76
$ws6.Range("B11").Value = 77
$ws6.Range("B12").Value = 78
$ws6.Range("B14").Value = 81
$ws6.Range("B15").Value = 82
$ws6.Range("B16").Value = 83
$ws6.Range("B18").Value = 86
$ws6.Range("B20").Value = 89
$ws6.Range("B22").Value = 92
$ws6.Range("B23").Value = 94
$ws6.Range("B24").Value = 96
$ws6.Range("B25").Value = 98
$ws6.Range("B28").Value = 108
$ws6.Range("B29").Value = 113

# ---- Sheet "8.0" (sheet7): targeted SS corrections ----
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("B2").Value  = 62
$ws7.Range("B3").Value  = 64
$ws7.Range("B4").Value  = 65
$ws7.Range("B5").Value  = 67
$ws7.Range("B6").Value  = 68
$ws7.Range("B7").Value  = 69
$ws7.Range("B8").Value  = 70
$ws7.Range("B9").Value  = 71
$ws7.Range("B10").Value = 73
$ws7.Range("B11").Value = 74
$ws7.Range("B12").Value = 75
$ws7.Range("B13").Value = 76
$ws7.Range("B14").Value = 78
$ws7.Range("B15").Value = 79
$ws7.Range("B16").Value = 80
$ws7.Range("B17").Value = 81
$ws7.Range("B18").Value = 83
$ws7.Range("B19").Value = 84
$ws7.Range("B20").Value = 86
$ws7.Range("B22").Value = 89
$ws7.Range("B23").Value = 91
$ws7.Range("B25").Value = 96
$ws7.Range("B27").Value = 105
$ws7.Range("B28").Value = 111

$ws1.Select()
